# Add two new "longname" entries (Transport Canada aircraft) on the
# "Sheet1" worksheet (the active/types-free data sheet, tabSelected=1),
# style them with an Arial/12pt/#212529 font, bump the row heights for
# rows 5-6, and leave the selection on F6 (matching the authored edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "longname" column values for the Transport Canada Dash 7 / Dash 8 rows
$ws.Range("F5").Value = "Transport Canada Dash 7 - CGCFR"
$ws.Range("F6").Value = "Transport Canada Dash 8 - CGCFJ"

# Apply the custom font (Arial 12, font color #212529) to F5 first...
$f5Chars = $ws.Range("F5").Characters()
$f5Chars.Font.Name = "Arial"
$f5Chars.Font.Size = 12
$f5Chars.Font.Color = 2696481   # RGB(33,37,41) == 0x212529 -> BGR long

# ...then copy just the formatting over to F6 so both cells end up sharing
# the same cell style instead of building up a second, independent one.
$ws.Range("F5").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Both edited rows grew slightly taller in the authored workbook
$ws.Rows("5:6").RowHeight = 15.75

# Match the final cursor/selection position left behind in the workbook
$ws.Range("F6").Select()
